# Burn Down Chart Actualizado
# Extend the "Trabajo restante" (remaining work) formulas/values in column C/D
# down through row 16, mirroring a fill-down of the existing C11 formula
# pattern (C = previous C - current D) and entering the new "Horas
# restantes" values for the last days of the sprint.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Re-enter C11 as an explicit (non-shared) formula, then fill the same
# remaining-work pattern down through row 16.
$ws.Range("C11").Formula = "=C10-D11"
$ws.Range("C12:C16").Formula = "=C11-D12"

# New daily remaining-hours entries for rows 12-16.
$ws.Range("D12").Value = 0
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("D16").Value = 10

# Rename the built-in "Bueno" cell style to "Buena".
$wb.Styles.Item("Bueno").Name = "Buena"

# Restore the active selection to the last edited cell.
$ws.Range("D16").Select()
